$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.795.65"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").Value = "2.416.85"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.83"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.54"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.84"
$ws.Range("E10").Value = "  +4.72%  "
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.74"
$ws.Range("E13").Value = "  +3.93%  "
$ws.Range("D14").Value = "2.847.99"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "59.640.59"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "2.382.96"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.41"
$ws.Range("E18").Value = "  +6.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.42"
$ws.Range("E19").Value = "  +3.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "333.49"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").Value = "  +2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.58"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.169"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.55"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Value = "0.0₃0785"
$ws.Range("E28").Value = "  +6.72%  "
$ws.Range("E29").Value = "  +3.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.81"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.70"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.02"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.32"
$ws.Range("E35").Value = "  +5.60%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.11"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("E40").Value = "  +11.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "309.97"
$ws.Range("E41").Value = "  +5.35%  "
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "143.33"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0964"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("E45").Value = "  +4.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.412"
$ws.Range("E46").Value = "  +6.95%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.572"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.62"
$ws.Range("E51").Value = "  +4.72%  "
